# Word COM-interop script applying the "Software test plan.docx" revision:
#  1) Split the run ending "...much better simulation..." so a _GoBack
#     bookmark sits between "much better " and "simulation which...".
#     (Word keeps a single _GoBack bookmark, so re-adding it here also
#     removes it from its old location - the blank paragraph under
#     "Overall results".)
#  2) Reword the "All previous posted team projects..." reference intro.
#  3) Mark "Ahamd" as a spell-check exception (proofErr spellStart/spellEnd)
#     and append " - Documentation" after "Ahamd Hawkins Probe".
#  4) Remove the "Firodiya Liu Kaul" / "Huang Anh Vedam" / "Li Pan Wang" /
#     "Mandal Parulker" / "Nowak Wilde" bullet paragraphs.
#  5) Append " - Crowd/UAV Spawn & Screen Ratio" after "Wang Zhou".

$d = $word.ActiveDocument

function Get-ParagraphText($p) {
    return $p.Range.Text.TrimEnd([char]13, [char]7)
}

function Find-ParagraphByText($doc, $text) {
    $n = $doc.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ((Get-ParagraphText $p) -eq $text) {
            return $p
        }
    }
    return $null
}

# Replace a whole paragraph's run content (everything except the trailing
# paragraph mark) with freshly authored run/proofErr XML. Clearing the
# range down to the paragraph's true start first is required - InsertXML
# only splices its <w:p> runs into the *existing* paragraph when the
# target range sits exactly at that paragraph's start; anywhere else it
# creates a new sibling paragraph instead of merging in place.
function Set-ParagraphRunXml($doc, $para, $innerXml) {
    $full = $doc.Range($para.Range.Start, $para.Range.End - 1)
    $full.Text = ""
    $frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $full.InsertXML($frag)
}

$szRpr = '<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'

# --- 1) Split "...much better / simulation..." run and place _GoBack ---

$splitAfter = "This effect provides for what we believe is a much better "
$fullText = $d.Content.Text
$splitIdx = $fullText.IndexOf($splitAfter)
if ($splitIdx -lt 0) {
    throw "Could not locate crowd-changes paragraph text to split"
}
$splitPoint = $splitIdx + $splitAfter.Length
$splitRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $splitRange) | Out-Null

# --- 2) Reword the references intro paragraph ---

$d.Content.Find.Execute(
    "All previous posted team projects contributed to some extent to this revision.  The groups whose projects we evaluated are recognized (listed as the group evaluations listed them):",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "We evaluated all the previous class projects and took what we though was relevant information from the following two teams.  ",
    2
) | Out-Null

# --- 3) "Ahamd Hawkins Probe" -> spell-flagged "Ahamd" + " Hawkins Probe" + " - Documentation" ---

$probePara = Find-ParagraphByText $d "Ahamd Hawkins Probe"
if ($null -eq $probePara) {
    throw "Could not locate 'Ahamd Hawkins Probe' paragraph"
}
$probeXml = '<w:proofErr w:type="spellStart"/>' +
    '<w:r>' + $szRpr + '<w:t>Ahamd</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r>' + $szRpr + '<w:t xml:space="preserve"> Hawkins Probe</w:t></w:r>' +
    '<w:r>' + $szRpr + '<w:t xml:space="preserve"> ' + [char]0x2013 + ' Documentation</w:t></w:r>'
Set-ParagraphRunXml $d $probePara $probeXml

# --- 4) Delete the five intervening bullet paragraphs ---

$firstDel = Find-ParagraphByText $d "Firodiya Liu Kaul"
$lastDel = Find-ParagraphByText $d "Nowak Wilde"
if (($null -eq $firstDel) -or ($null -eq $lastDel)) {
    throw "Could not locate bullet paragraphs to delete"
}
$delRange = $d.Range($firstDel.Range.Start, $lastDel.Range.End)
$delRange.Delete() | Out-Null

# --- 5) Append " - Crowd/UAV Spawn & Screen Ratio" after "Wang Zhou" ---

$zhouPara = Find-ParagraphByText $d "Wang Zhou"
if ($null -eq $zhouPara) {
    throw "Could not locate 'Wang Zhou' paragraph"
}
$zhouXml = '<w:r>' + $szRpr + '<w:t>Wang Zhou</w:t></w:r>' +
    '<w:r>' + $szRpr + '<w:t xml:space="preserve"> ' + [char]0x2013 + ' Crowd/UAV Spawn &amp; Screen Ratio</w:t></w:r>'
Set-ParagraphRunXml $d $zhouPara $zhouXml

Write-Output "edit complete"
